$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Narrow the tab ratio gap between open sheet tabs and the horizontal
# scrollbar slightly (workbookView tabRatio 989 -> 990).
$excel.ActiveWindow.TabRatio = 0.99

# A student batch label was typed into E2 ("2018-19_ODD"), which also
# becomes a new shared string.
$ws.Range("E2").Value = "2018-19_ODD"

# The active cell/selection moved from G2 to F2.
[void]$ws.Range("F2").Select()

# Columns A, B and D were nudged slightly wider (autofit-style tweak).
$ws.Columns.Item(1).ColumnWidth = 18.126180836707167
$ws.Columns.Item(2).ColumnWidth = 24.340755735492568
$ws.Columns.Item(4).ColumnWidth = 12.555330634277967
